$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.659.97"
$ws.Range("D3").Value = "'1.591.33"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'211.16"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").Value = "'0.0835"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "'1.814.25"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").Value = "'1.588.53"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("D16").Value = "'64.73"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'26.663.38"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'208.01"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "'147.17"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").Value = "'0.664"
$ws.Range("E33").Value = "  +22.80%  "
$ws.Range("D34").Value = "'1.325.07"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'0.827"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("D42").Value = "'0.787"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'1.727.13"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'0.837"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  -0.70%  "
